$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update progress values from 0 to 1 for the tracked rows
$rows = @(5,6,7,9,10,11,12,13,15,16,17,18)
foreach ($r in $rows) {
    $ws.Range("C$r").Value = 1
}

# Fix the AVERAGE formula (remove trailing comma / empty argument)
$ws.Range("F4").Formula = "=AVERAGE(C5:C7,C9:C13,C15:C18)"

# Update the active view/selection to match the latest edit state
# (scrolled so row 10 is at the top, with B11 as the active selected cell)
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B11").Select()

$wb.Save()
